# Auto-generated edit script: updates market-price / profit columns (H-N)
# on the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to reflect refreshed
# marketboard data, per the scheduled runner's update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1232.7273
$ws.Range("I6").Value = 1232.7273
$ws.Range("K6").Value = 3698.1819
$ws.Range("M6").Value = -3586.1819

$ws.Range("H12").Value = 133.21428
$ws.Range("I12").Value = 97.083336
$ws.Range("J12").Value = 350
$ws.Range("K12").Value = 97.083336
$ws.Range("L12").Value = 350
$ws.Range("M12").Value = 72.916664
$ws.Range("N12").Value = -690

$ws.Range("H21").Value = 36603.8
$ws.Range("I21").Value = 25754.75
$ws.Range("K21").Value = 25754.75
$ws.Range("M21").Value = -25286.75

$ws.Range("H23").Value = 36603.8
$ws.Range("I23").Value = 25754.75
$ws.Range("K23").Value = 25754.75
$ws.Range("M23").Value = -25520.75

$ws.Range("H29").Value = 266
$ws.Range("I29").Value = 82.5
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 247.5
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = 33.5
$ws.Range("N29").Value = -3562

$ws.Range("H38").Value = 2844.9546
$ws.Range("I38").Value = 113.42857
$ws.Range("J38").Value = 4119.6665
$ws.Range("K38").Value = 340.28571
$ws.Range("L38").Value = 12358.9995
$ws.Range("M38").Value = 31.71429000000001
$ws.Range("N38").Value = -13102.9995

$ws.Range("H43").Value = 126923.5
$ws.Range("I43").Value = 200300.2
$ws.Range("J43").Value = 4629
$ws.Range("K43").Value = 200300.2
$ws.Range("L43").Value = 4629
$ws.Range("M43").Value = -200231.2
$ws.Range("N43").Value = -4767

$ws.Range("H58").Value = 1118.1578
$ws.Range("I58").Value = 206.5
$ws.Range("J58").Value = 2131.111
$ws.Range("K58").Value = 619.5
$ws.Range("L58").Value = 6393.333
$ws.Range("M58").Value = -469.5
$ws.Range("N58").Value = -6693.333

$ws.Range("H87").Value = 13020.853
$ws.Range("J87").Value = 13020.853
$ws.Range("L87").Value = 13020.853
$ws.Range("N87").Value = -15516.853

$ws.Range("H90").Value = 13020.853
$ws.Range("J90").Value = 13020.853
$ws.Range("L90").Value = 39062.55899999999
$ws.Range("N90").Value = -51542.55899999999

$ws.Range("H137").Value = 14913385
$ws.Range("I137").Value = 1050.7037
$ws.Range("J137").Value = 51516388
$ws.Range("K137").Value = 3152.1111
$ws.Range("L137").Value = 154549164
$ws.Range("M137").Value = -602.1111000000001
$ws.Range("N137").Value = -154554264


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 51112736
$ws.Range("I74").Value = 41667388
$ws.Range("J74").Value = 88894130
$ws.Range("K74").Value = 41667388
$ws.Range("L74").Value = 88894130
$ws.Range("M74").Value = -41666514
$ws.Range("N74").Value = -88895878

$ws.Range("H77").Value = 51112736
$ws.Range("I77").Value = 41667388
$ws.Range("J77").Value = 88894130
$ws.Range("K77").Value = 208336940
$ws.Range("L77").Value = 444470650
$ws.Range("M77").Value = -208332572
$ws.Range("N77").Value = -444479386

$ws.Range("H97").Value = 620.6
$ws.Range("I97").Value = 619
$ws.Range("K97").Value = 619
$ws.Range("M97").Value = -123

$ws.Range("H132").Value = 14623801
$ws.Range("I132").Value = 16670849
$ws.Range("J132").Value = 6947369
$ws.Range("K132").Value = 50012547
$ws.Range("L132").Value = 20842107
$ws.Range("M132").Value = -50010017
$ws.Range("N132").Value = -20847167


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 22112
$ws.Range("J81").Value = 22112
$ws.Range("L81").Value = 22112
$ws.Range("N81").Value = -24234

$ws.Range("H84").Value = 22112
$ws.Range("J84").Value = 22112
$ws.Range("L84").Value = 66336
$ws.Range("N84").Value = -76944

$ws.Range("H94").Value = 1620.2778
$ws.Range("I94").Value = 1079.0625
$ws.Range("K94").Value = 1079.0625
$ws.Range("M94").Value = -628.0625

$ws.Range("H99").Value = 1239.762
$ws.Range("I99").Value = 1068.7142
$ws.Range("J99").Value = 1325.2858
$ws.Range("K99").Value = 1068.7142
$ws.Range("L99").Value = 1325.2858
$ws.Range("M99").Value = 429.2858000000001
$ws.Range("N99").Value = -4321.2858


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1692656.8
$ws.Range("I31").Value = 1158.7407
$ws.Range("J31").Value = 6259701.5
$ws.Range("K31").Value = 1158.7407
$ws.Range("L31").Value = 6259701.5
$ws.Range("M31").Value = -863.7407000000001
$ws.Range("N31").Value = -6260291.5

$ws.Range("H34").Value = 1692656.8
$ws.Range("I34").Value = 1158.7407
$ws.Range("J34").Value = 6259701.5
$ws.Range("K34").Value = 1158.7407
$ws.Range("L34").Value = 6259701.5
$ws.Range("M34").Value = -956.7407000000001
$ws.Range("N34").Value = -6260105.5

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0

$ws.Range("H132").Value = 1077.5
$ws.Range("I132").Value = 786.5405
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 2359.6215
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = 170.3785000000003
$ws.Range("N132").Value = -19058

$ws.Range("H134").Value = 741647.4399999999
$ws.Range("I134").Value = 897.24
$ws.Range("K134").Value = 2691.72
$ws.Range("M134").Value = -156.7200000000003

$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1020
$ws.Range("I17").Value = 275
$ws.Range("J17").Value = 4000
$ws.Range("K17").Value = 825
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = -656
$ws.Range("N17").Value = -12338

$ws.Range("H34").Value = 2469.1538
$ws.Range("I34").Value = 300
$ws.Range("J34").Value = 2649.9167
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 7949.750100000001
$ws.Range("M34").Value = -816
$ws.Range("N34").Value = -8117.750100000001

$ws.Range("H39").Value = 2056.8
$ws.Range("I39").Value = 446.66666
$ws.Range("J39").Value = 2459.3333
$ws.Range("K39").Value = 1339.99998
$ws.Range("L39").Value = 7377.999899999999
$ws.Range("M39").Value = -1045.99998
$ws.Range("N39").Value = -7965.999899999999

$ws.Range("H55").Value = 2350.5
$ws.Range("I55").Value = 6
$ws.Range("J55").Value = 2563.6365
$ws.Range("K55").Value = 18
$ws.Range("L55").Value = 7690.9095
$ws.Range("M55").Value = 159
$ws.Range("N55").Value = -8044.9095

$ws.Range("H92").Value = 1315502.8
$ws.Range("I92").Value = 169.63637
$ws.Range("K92").Value = 508.90911
$ws.Range("M92").Value = 739.0908899999999

$ws.Range("H113").Value = 1173.6
$ws.Range("I113").Value = 940.9167
$ws.Range("J113").Value = 1309.8049
$ws.Range("K113").Value = 2822.7501
$ws.Range("L113").Value = 3929.4147
$ws.Range("M113").Value = -652.7501000000002
$ws.Range("N113").Value = -8269.414700000001

$ws.Range("H134").Value = 2758.8
$ws.Range("I134").Value = 1948.75
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 5846.25
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = -776.25
$ws.Range("N134").Value = -28137


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14106402
$ws.Range("I132").Value = 13033319
$ws.Range("J132").Value = 18184120
$ws.Range("K132").Value = 39099957
$ws.Range("L132").Value = 54552360
$ws.Range("M132").Value = -39097427
$ws.Range("N132").Value = -54557420


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3985.7036
$ws.Range("I22").Value = 2514.2
$ws.Range("J22").Value = 4851.294
$ws.Range("K22").Value = 2514.2
$ws.Range("L22").Value = 4851.294
$ws.Range("M22").Value = -2219.2
$ws.Range("N22").Value = -5441.294

$ws.Range("H27").Value = 3985.7036
$ws.Range("I27").Value = 2514.2
$ws.Range("J27").Value = 4851.294
$ws.Range("K27").Value = 2514.2
$ws.Range("L27").Value = 4851.294
$ws.Range("M27").Value = -2407.2
$ws.Range("N27").Value = -5065.294

$ws.Range("H46").Value = 1137.9286
$ws.Range("I46").Value = 1123.5
$ws.Range("J46").Value = 1148.75
$ws.Range("K46").Value = 1123.5
$ws.Range("L46").Value = 1148.75
$ws.Range("M46").Value = -935.5
$ws.Range("N46").Value = -1524.75

$ws.Range("H55").Value = 35714670
$ws.Range("I55").Value = 250000000
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 250000000
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = -249999827
$ws.Range("N55").Value = -796

$ws.Range("H80").Value = 27999
$ws.Range("J80").Value = 27999
$ws.Range("L80").Value = 27999
$ws.Range("N80").Value = -30245

$ws.Range("H83").Value = 27999
$ws.Range("J83").Value = 27999
$ws.Range("L83").Value = 83997
$ws.Range("N83").Value = -95229

$ws.Range("H132").Value = 2507298
$ws.Range("I132").Value = 3040166.5
$ws.Range("K132").Value = 9120499.5
$ws.Range("M132").Value = -9117969.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1071.4651
$ws.Range("I136").Value = 369.82
$ws.Range("K136").Value = 1109.46
$ws.Range("M136").Value = 1440.54

